# Generate Report for Handback
# Update the "generated at" timestamps written into the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first data row.
# The same text is also used by de-de!H2 (Correspond Handoff Datetime), so
# both cells are updated to keep them sharing the same string value.
$wsOverview.Range("G2").Value = "2016-08-20 13:06:27"
$wsDeDe.Range("H2").Value = "2016-08-20 13:06:27"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for row 2.
$wsZhCn.Range("H2").Value = "2016-08-20 13:06:22"
$wsZhCn.Range("K2").Value = "2016-08-20 13:06:39"

# de-de sheet: Correspond Handback DateTime for row 2.
$wsDeDe.Range("K2").Value = "2016-08-20 13:06:45"
